$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 182 (pushes existing rows 182.. down to 184..)
$ws.Rows.Item(182).Insert()
$ws.Rows.Item(182).Insert()

# New row 182 data
$ws.Cells.Item(182, 1).Value = 3
$ws.Cells.Item(182, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 44455
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 100112021
$ws.Cells.Item(182, 7).Value = "Ají"
$ws.Cells.Item(182, 8).Value = "Americana (o)"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 25
$ws.Cells.Item(182, 11).Value = 73000
$ws.Cells.Item(182, 12).Value = 73000
$ws.Cells.Item(182, 13).Value = 73000
$ws.Cells.Item(182, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(182, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(182, 16).Value = 2920
$ws.Cells.Item(182, 17).Value = 25
$ws.Cells.Item(182, 18).Value = "Hortaliza"

# New row 183 data
$ws.Cells.Item(183, 1).Value = 3
$ws.Cells.Item(183, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(183, 3).Value = "Coquimbo"
$ws.Cells.Item(183, 4).Value = 44455
$ws.Cells.Item(183, 5).Value = 5
$ws.Cells.Item(183, 6).Value = 100112021
$ws.Cells.Item(183, 7).Value = "Ají"
$ws.Cells.Item(183, 8).Value = "Inferno"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 35
$ws.Cells.Item(183, 11).Value = 40000
$ws.Cells.Item(183, 12).Value = 40000
$ws.Cells.Item(183, 13).Value = 40000
$ws.Cells.Item(183, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(183, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(183, 16).Value = 2667
$ws.Cells.Item(183, 17).Value = 15
$ws.Cells.Item(183, 18).Value = "Hortaliza"
